$d = $word.ActiveDocument

# wdReplace constants used below: wdReplaceNone=0, wdReplaceOne=1, wdReplaceAll=2

# 1) Header table cell: "PHÒNG HÀNH CHÍNH" -> "PHÒNG ĐÀO TẠO"
#    This exact (mixed single-space, all-caps) string occurs only once in the
#    document, so a document-wide replace-all is safe and precise.
$d.Content.Find.Execute("PHÒNG HÀNH CHÍNH", $true, $false, $false, $false, $false, `
    $true, 1, $false, "PHÒNG ĐÀO TẠO", 2) | Out-Null

# 2) Date: "05/11/2024" -> "10/12/2024"
#    Also a unique, single occurrence in the document.
$d.Content.Find.Execute("05/11/2024", $true, $false, $false, $false, $false, `
    $true, 1, $false, "10/12/2024", 2) | Out-Null

# 3) "Phòng Hành chính" -> "Phòng Đào tạo"
#    This lower-case phrase appears 4 times in the document, but only the 2nd
#    occurrence ("Ngày ..., Phòng Hành chính nhận được Công văn sau:") and the
#    4th occurrence ("Ý kiến của Phòng Hành chính: ...") change; the 1st
#    ("Văn thư đơn vị Phòng Hành chính.") and 3rd ("CQ phát hành: Phòng Hành
#    chính") stay exactly as-is. So we walk the matches one at a time (case
#    sensitive, so the earlier all-caps heading is never touched here) and
#    replace only the 2nd/4th hit, using a duplicate Range for the actual
#    replace so the walking cursor itself is not disturbed.
$rng = $d.Content
$rng.Start = 0
$rng.End = $d.Content.End
$idx = 0
$guard = 0
while ($rng.Find.Execute("Phòng Hành chính", $true, $false, $false, $false, $false, `
        $true, 1, $false, "", 0)) {
    $idx = $idx + 1
    if ($idx -eq 2 -or $idx -eq 4) {
        $target = $d.Range($rng.Start, $rng.End)
        $target.Text = "Phòng Đào tạo"
    }
    $rng.Start = $rng.End
    $rng.End = $d.Content.End

    $guard = $guard + 1
    if ($guard -gt 20) { break }
}

Write-Output "Replacements complete (matched $idx occurrences of 'Phong Hanh chinh')"
